# BOOK STORE INSTALLATION GUIDE.docx — section 2 ("INSTALLATION PREREQUISITES")
# content refresh, plus small fixes further down in section 3.
#
# Semantic summary of the change (see diff):
#   1. Remove one of the two blank lines directly above the
#      "2. INSTALLATION PREREQUISITES" heading (was two blank NoSpacing
#      paragraphs, now only one).
#   2. Intro paragraph: "hardware and software requirements" ->
#      "hardware, database and software requirements".
#   3. Hardware list: "Processor:" -> "Processor: Core 2 Duo",
#      "Memory:" -> "Memory: 4GB", and the "Disk:" line is removed.
#   4. "MySQL/MariaDB " -> "MySQL/MariaDB 10.4".
#   5. Remove the blank line above the software list and replace the
#      three old software lines with "Apache 2.4", "PHP 7.4",
#      "Windows Server 2016".
#   6. "Run the installer program on the server." -> "Copy program
#      files and folders on the server."
#   7. Remove the blank line directly under the "Post-Installation
#      Tasks" heading.

$d = $word.ActiveDocument

function TextOf($para) {
    # Paragraph.Range.Text always carries the trailing paragraph mark
    # (chr 13); strip it so comparisons are against the visible text.
    return $para.Range.Text.TrimEnd([char]13)
}

function DeleteParagraph($para) {
    $d.Range($para.Range.Start, $para.Range.End).Delete()
}

# ---------------------------------------------------------------------
# 1. Collapse the double blank line above "2. INSTALLATION PREREQUISITES"
#    into a single blank line.
# ---------------------------------------------------------------------
$installPrereqHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 1" -and (TextOf $p) -like "*2. INSTALLATION PREREQUISITES*") {
        $installPrereqHeading = $p
    }
}
if ($installPrereqHeading -ne $null) {
    $blankAbove = $installPrereqHeading.Previous()
    if ($blankAbove -ne $null -and (TextOf $blankAbove) -eq "") {
        DeleteParagraph $blankAbove
    }
}

# ---------------------------------------------------------------------
# 2. Intro paragraph wording.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "This chapter lists all hardware and software requirements for the installation of ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This chapter lists all hardware, database and software requirements for the installation of ",
    2) | Out-Null

# ---------------------------------------------------------------------
# 3. Hardware requirements list.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Processor:", $true, $false, $false, $false, $false, $true, 1, $false,
    "Processor: Core 2 Duo", 2) | Out-Null

$d.Content.Find.Execute(
    "Memory:", $true, $false, $false, $false, $false, $true, 1, $false,
    "Memory: 4GB", 2) | Out-Null

$diskParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ((TextOf $p) -eq "Disk:") {
        $diskParagraph = $p
    }
}
if ($diskParagraph -ne $null) {
    DeleteParagraph $diskParagraph
}

# ---------------------------------------------------------------------
# 4. Database prerequisites: MySQL/MariaDB version.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "MySQL/MariaDB ", $true, $false, $false, $false, $false, $true, 1, $false,
    "MySQL/MariaDB 10.4", 2) | Out-Null

# ---------------------------------------------------------------------
# 5. Software requirements list: drop the blank separator line, then
#    swap the three sample products for the new ones.
# ---------------------------------------------------------------------
$identifySoftware = $null
foreach ($p in $d.Paragraphs) {
    if ((TextOf $p) -eq "Identify software prerequisites, such as:") {
        $identifySoftware = $p
    }
}
if ($identifySoftware -ne $null) {
    $blankBelow = $identifySoftware.Next()
    if ($blankBelow -ne $null -and (TextOf $blankBelow) -eq "") {
        DeleteParagraph $blankBelow
    }
}

$d.Content.Find.Execute(
    "Microsoft® Windows Server™ 2003", $true, $false, $false, $false, $false, $true, 1, $false,
    "Apache 2.4", 2) | Out-Null

$d.Content.Find.Execute(
    "Microsoft .NET Framework 3.0", $true, $false, $false, $false, $false, $true, 1, $false,
    "PHP 7.4", 2) | Out-Null

$d.Content.Find.Execute(
    "Microsoft SQL Server™ 2010", $true, $false, $false, $false, $false, $true, 1, $false,
    "Windows Server 2016", 2) | Out-Null

# ---------------------------------------------------------------------
# 6. Installation step wording.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Run the installer program on the server.", $true, $false, $false, $false, $false, $true, 1, $false,
    "Copy program files and folders on the server.", 2) | Out-Null

# ---------------------------------------------------------------------
# 7. Remove the blank line under "Post-Installation Tasks".
# ---------------------------------------------------------------------
$postInstallHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 2" -and (TextOf $p) -like "*Post-Installation Tasks*") {
        $postInstallHeading = $p
    }
}
if ($postInstallHeading -ne $null) {
    $blankUnder = $postInstallHeading.Next()
    if ($blankUnder -ne $null -and (TextOf $blankUnder) -eq "") {
        DeleteParagraph $blankUnder
    }
}

Write-Output "done"
